$d = $word.ActiveDocument
$w_ns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------------
# This document review thread (the "nicotine material" / BIOBANK_0000048
# entry) gets three changes, all scoped to the small cluster of paragraphs
# right after the "Shall we change to 'material entity' ..." question:
#
#   1. A new reviewer comment run is appended to that question paragraph.
#   2. The big red "DONE" paragraph right below it is reworded to
#      "No change made" (the nicotine-material term change is being
#      reversed per the commit message), and the "_GoBack" bookmark that
#      used to sit two paragraphs further down is pulled up into this
#      paragraph.
#   3. The now-bookmark-less paragraph becomes a plain empty paragraph,
#      matching its sibling right above it.
# ---------------------------------------------------------------------------

# Locate the paragraphs by their known, stable absolute indices in this
# document and sanity-check their text before touching anything, so the
# script fails loudly instead of silently editing the wrong spot.
$pQuestion = $d.Paragraphs.Item(75)
$pDone     = $d.Paragraphs.Item(76)
$pBlank    = $d.Paragraphs.Item(77)
$pGoBack   = $d.Paragraphs.Item(78)

$cr = [char]13
$questionText = $pQuestion.Range.Text.TrimEnd($cr)
if ($questionText -notmatch "Shall we change to 'material entity'") {
    throw "Expected paragraph 75 to be the 'Shall we change...' question, got: $questionText"
}
if ($pDone.Range.Text.TrimEnd($cr) -ne "DONE") {
    throw "Expected paragraph 76 to read DONE, got: $($pDone.Range.Text)"
}
if ($pBlank.Range.Text.TrimEnd($cr) -ne "") {
    throw "Expected paragraph 77 to be blank, got: $($pBlank.Range.Text)"
}
if ($pGoBack.Range.Text.TrimEnd($cr) -ne "") {
    throw "Expected paragraph 78 to be blank (bookmark only), got: $($pGoBack.Range.Text)"
}

# 1) Append the new reviewer comment as its own run on the question
#    paragraph, preserving the existing tab + text run untouched.
$newComment = " Based on textual definition, should be " + [char]0x2018 + "processed material" + [char]0x2019
$questionXml = "<w:p xmlns:w='$w_ns'>" +
    "<w:pPr><w:spacing w:after=`"0`"/></w:pPr>" +
    "<w:r><w:tab/><w:t>Shall we change to 'material entity' and ('has part' some nicotine) to make it more general?</w:t></w:r>" +
    "<w:r><w:t xml:space=`"preserve`">$newComment</w:t></w:r>" +
    "</w:p>"
$pQuestion.Range.InsertXML($questionXml)

# 2) Remove the "_GoBack" bookmark from its own paragraph first (so the
#    id=0 slot is free and the bookmark we add next keeps that same id),
#    turning that paragraph into a plain empty one.
$goBackXml = "<w:p xmlns:w='$w_ns'><w:pPr><w:spacing w:after=`"0`"/></w:pPr></w:p>"
$pGoBack.Range.InsertXML($goBackXml)

# 3) Reword DONE -> No change made, keeping the bold/red run formatting,
#    and move the "_GoBack" bookmark into this paragraph (right after the
#    run, same as the target markup).
$doneXml = "<w:p xmlns:w='$w_ns'>" +
    "<w:pPr><w:spacing w:after=`"0`"/><w:rPr><w:b/><w:color w:val=`"FF0000`"/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:b/><w:color w:val=`"FF0000`"/></w:rPr><w:t>No change made</w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
    "</w:p>"
$pDone.Range.InsertXML($doneXml)

Write-Output "edit applied"
